# Updated cryptos list values (Price / Volume(1h)) per the target diff.
# Cells hold plain numeric-looking TEXT (not numbers), so each cell is
# temporarily forced to Text format before the write, then ClearFormats()
# is used to drop the temporary style again (matches the original "no
# explicit style" cells) while Excel keeps the stored value as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.716.37'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("E2").ClearFormats()
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.848.48'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E3").ClearFormats()
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -2.76%  '
$ws.Range("E4").ClearFormats()
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.41'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("E5").ClearFormats()
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.80%  '
$ws.Range("E6").ClearFormats()
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4316'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.50%  '
$ws.Range("E7").ClearFormats()
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3745'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("E8").ClearFormats()
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("E9").ClearFormats()
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8809'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E10").ClearFormats()
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E11").ClearFormats()
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.847.46'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("E12").ClearFormats()
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.737'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E13").ClearFormats()
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.451'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("E14").ClearFormats()
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07145'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("E15").ClearFormats()
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.04'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.94%  '
$ws.Range("E16").ClearFormats()
# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.78%  '
$ws.Range("E17").ClearFormats()
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008997'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("E18").ClearFormats()
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.009'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("E19").ClearFormats()
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.48'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("E20").ClearFormats()
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.714.97'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E21").ClearFormats()
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.250'
$ws.Range("D22").ClearFormats()
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.074.96'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("E24").ClearFormats()
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("E25").ClearFormats()
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.56'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("E26").ClearFormats()
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.61'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("E27").ClearFormats()
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.142'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.68%  '
$ws.Range("E28").ClearFormats()
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.396'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("E29").ClearFormats()
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.34'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.04%  '
$ws.Range("E30").ClearFormats()
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08938'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.59%  '
$ws.Range("E31").ClearFormats()
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.232'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.33%  '
$ws.Range("E32").ClearFormats()
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7793'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("E33").ClearFormats()
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.567'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("E34").ClearFormats()
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.916'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.48%  '
$ws.Range("E35").ClearFormats()
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.011'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.69%  '
$ws.Range("E36").ClearFormats()
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.141'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.48%  '
$ws.Range("E37").ClearFormats()
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05337'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("E38").ClearFormats()
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E39").ClearFormats()
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.266'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.35%  '
$ws.Range("E40").ClearFormats()
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.874'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("E41").ClearFormats()
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5163'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("E42").ClearFormats()
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1678'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("E43").ClearFormats()
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.922'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("E44").ClearFormats()
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.64'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("E45").ClearFormats()
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.63'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("E46").ClearFormats()
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4735'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("E47").ClearFormats()
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06509'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("E48").ClearFormats()
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.699'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("E49").ClearFormats()
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.012'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.68%  '
$ws.Range("E50").ClearFormats()
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.883'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.80%  '
$ws.Range("E51").ClearFormats()
